$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update BG column (Number of machinery) values
$ws.Range("BG2").Value = 34
$ws.Range("BG3").Value = 38
$ws.Range("BG5").Value = 1
$ws.Range("BG6").Value = 25
$ws.Range("BG7").Value = 100
$ws.Range("BG8").Value = 36
$ws.Range("BG9").Value = 17
$ws.Range("BG19").Value = 29
$ws.Range("BG20").Value = 88
$ws.Range("BG24").Value = 25
$ws.Range("BG25").Value = 25
$ws.Range("BG27").Value = 67
$ws.Range("BG31").Value = 62
$ws.Range("BG32").Value = 9
$ws.Range("BG33").Value = 18
$ws.Range("BG34").Value = 26
$ws.Range("BG36").Value = 1
$ws.Range("BG37").Value = 72
$ws.Range("BG38").Value = 60
$ws.Range("BG39").Value = 49
$ws.Range("BG40").Value = 49
$ws.Range("BG41").Value = 36
$ws.Range("BG42").Value = 40
$ws.Range("BG43").Value = 24
$ws.Range("BG45").Value = 3
$ws.Range("BG46").Value = 1
$ws.Range("BG47").Value = 206
$ws.Range("BG48").Value = 44
$ws.Range("BG49").Value = 4
$ws.Range("BG50").Value = 82
$ws.Range("BG51").Value = 101
$ws.Range("BG52").Value = 20
$ws.Range("BG53").Value = 36
$ws.Range("BG55").Value = 25
$ws.Range("BG56").Value = 12
$ws.Range("BG60").Value = 99
$ws.Range("BG61").Value = 25
$ws.Range("BG63").Value = 32
$ws.Range("BG64").Value = 46
$ws.Range("BG65").Value = 8
$ws.Range("BG67").Value = 56
$ws.Range("BG68").Value = 30
$ws.Range("BG69").Value = 17
$ws.Range("BG70").Value = 79
$ws.Range("BG71").Value = 24
$ws.Range("BG73").Value = 118
$ws.Range("BG75").Value = 10
$ws.Range("BG77").Value = 119
$ws.Range("BG80").Value = 31
$ws.Range("BG81").Value = 239
$ws.Range("BG82").Value = 82
$ws.Range("BG87").Value = 68
$ws.Range("BG88").Value = 151
$ws.Range("BG89").Value = 120
$ws.Range("BG90").Value = 93
$ws.Range("BG91").Value = 11
$ws.Range("BG92").Value = 72
$ws.Range("BG93").Value = 120
$ws.Range("BG94").Value = 93
$ws.Range("BG95").Value = 57
$ws.Range("BG96").Value = 21
$ws.Range("BG97").Value = 36
$ws.Range("BG98").Value = 10
$ws.Range("BG102").Value = 43
$ws.Range("BG104").Value = 66
$ws.Range("BG105").Value = 100
$ws.Range("BG106").Value = 41
$ws.Range("BG107").Value = 34
$ws.Range("BG108").Value = 123
$ws.Range("BG109").Value = 245
$ws.Range("BG110").Value = 70
$ws.Range("BG111").Value = 79
$ws.Range("BG112").Value = 63
$ws.Range("BG113").Value = 85
$ws.Range("BG114").Value = 35
$ws.Range("BG115").Value = 14
$ws.Range("BG116").Value = 31
$ws.Range("BG117").Value = 52
$ws.Range("BG118").Value = 39
$ws.Range("BG119").Value = 36
$ws.Range("BG121").Value = 85
$ws.Range("BG122").Value = 79
$ws.Range("BG123").Value = 82
$ws.Range("BG124").Value = 92
$ws.Range("BG125").Value = 76
$ws.Range("BG128").Value = 2
$ws.Range("BG131").Value = 55
$ws.Range("BG132").Value = 86
$ws.Range("BG133").Value = 96
$ws.Range("BG134").Value = 552
$ws.Range("BG135").Value = 5
$ws.Range("BG136").Value = 36
$ws.Range("BG137").Value = 4
$ws.Range("BG138").Value = 357
$ws.Range("BG139").Value = 465
$ws.Range("BG140").Value = 109
$ws.Range("BG142").Value = 1
$ws.Range("BG155").Value = 59
$ws.Range("BG156").Value = 27
$ws.Range("BG157").Value = 40
$ws.Range("BG158").Value = 34
$ws.Range("BG159").Value = 225
$ws.Range("BG163").Value = 21
$ws.Range("BG164").Value = 44
$ws.Range("BG165").Value = 17
$ws.Range("BG166").Value = 73
$ws.Range("BG167").Value = 47
$ws.Range("BG168").Value = 6
$ws.Range("BG170").Value = 43
$ws.Range("BG173").Value = 11
$ws.Range("BG174").Value = 181
$ws.Range("BG175").Value = 269
$ws.Range("BG177").Value = 168
$ws.Range("BG178").Value = 164
$ws.Range("BG179").Value = 1
$ws.Range("BG180").Value = 33
$ws.Range("BG182").Value = 149
$ws.Range("BG183").Value = 98
$ws.Range("BG184").Value = 66
$ws.Range("BG185").Value = 131
$ws.Range("BG188").Value = 1
$ws.Range("BG189").Value = 5
$ws.Range("BG192").Value = 51
$ws.Range("BG194").Value = 70
$ws.Range("BG195").Value = 4
$ws.Range("BG196").Value = 88
$ws.Range("BG197").Value = 19
$ws.Range("BG202").Value = 15
$ws.Range("BG203").Value = 103
$ws.Range("BG204").Value = 29
$ws.Range("BG205").Value = 3
$ws.Range("BG206").Value = 39
$ws.Range("BG210").Value = 27
$ws.Range("BG211").Value = 117
$ws.Range("BG212").Value = 41
$ws.Range("BG214").Value = 77
$ws.Range("BG215").Value = 79
$ws.Range("BG216").Value = 71
$ws.Range("BG217").Value = 3
$ws.Range("BG218").Value = 52
$ws.Range("BG219").Value = 64
$ws.Range("BG220").Value = 66
$ws.Range("BG221").Value = 97
$ws.Range("BG222").Value = 108
$ws.Range("BG223").Value = 36
$ws.Range("BG225").Value = 18
$ws.Range("BG226").Value = 9
$ws.Range("BG227").Value = 10
$ws.Range("BG229").Value = 39
$ws.Range("BG230").Value = 56
$ws.Range("BG231").Value = 4
$ws.Range("BG232").Value = 154
$ws.Range("BG233").Value = 66
$ws.Range("BG236").Value = 30
$ws.Range("BG238").Value = 307
$ws.Range("BG240").Value = 51
$ws.Range("BG241").Value = 54
$ws.Range("BG242").Value = 96
$ws.Range("BG245").Value = 25
$ws.Range("BG246").Value = -1
$ws.Range("BG249").Value = 40
$ws.Range("BG251").Value = 66
$ws.Range("BG252").Value = 202
$ws.Range("BG253").Value = 101
$ws.Range("BG255").Value = 51
$ws.Range("BG256").Value = 65
$ws.Range("BG257").Value = 77
$ws.Range("BG261").Value = 36
$ws.Range("BG262").Value = 47
$ws.Range("BG263").Value = 131
$ws.Range("BG264").Value = 85
$ws.Range("BG266").Value = 301
$ws.Range("BG267").Value = 87
$ws.Range("BG268").Value = 16
$ws.Range("BG272").Value = 28
$ws.Range("BG273").Value = 61
$ws.Range("BG274").Value = 43
$ws.Range("BG275").Value = 11
$ws.Range("BG276").Value = 114
$ws.Range("BG280").Value = 15
$ws.Range("BG281").Value = 14
$ws.Range("BG283").Value = 33
$ws.Range("BG284").Value = 2
$ws.Range("BG286").Value = 115
$ws.Range("BG287").Value = 14
$ws.Range("BG290").Value = 272
$ws.Range("BG292").Value = 60
$ws.Range("BG297").Value = 8
$ws.Range("BG298").Value = 26
$ws.Range("BG299").Value = 43
$ws.Range("BG302").Value = 16
$ws.Range("BG303").Value = 194
$ws.Range("BG304").Value = 59
$ws.Range("BG309").Value = 317
$ws.Range("BG310").Value = 87
$ws.Range("BG311").Value = 150
$ws.Range("BG312").Value = 74
$ws.Range("BG313").Value = 16
$ws.Range("BG315").Value = 16
$ws.Range("BG316").Value = 39
$ws.Range("BG317").Value = 121
$ws.Range("BG318").Value = 1
$ws.Range("BG319").Value = 2
$ws.Range("BG320").Value = 6
$ws.Range("BG322").Value = 8
$ws.Range("BG325").Value = 2
$ws.Range("BG328").Value = 15
$ws.Range("BG330").Value = 5
$ws.Range("BG331").Value = 40
$ws.Range("BG332").Value = 34
$ws.Range("BG333").Value = 29
$ws.Range("BG334").Value = 19
$ws.Range("BG335").Value = 20
$ws.Range("BG338").Value = 42
$ws.Range("BG339").Value = 20
$ws.Range("BG340").Value = 10
$ws.Range("BG342").Value = 16
$ws.Range("BG343").Value = 49
$ws.Range("BG346").Value = 84
$ws.Range("BG347").Value = 101
$ws.Range("BG348").Value = 53
$ws.Range("BG349").Value = 35
$ws.Range("BG350").Value = 58
$ws.Range("BG351").Value = 61
$ws.Range("BG354").Value = 125
$ws.Range("BG355").Value = 26
$ws.Range("BG357").Value = 1
$ws.Range("BG358").Value = 13
$ws.Range("BG359").Value = 47
$ws.Range("BG360").Value = 123
$ws.Range("BG361").Value = 9
$ws.Range("BG362").Value = 3
$ws.Range("BG363").Value = 81
$ws.Range("BG364").Value = 35
$ws.Range("BG368").Value = 14
$ws.Range("BG369").Value = 95
$ws.Range("BG370").Value = 57
$ws.Range("BG371").Value = 29
$ws.Range("BG372").Value = 13
$ws.Range("BG374").Value = 20
$ws.Range("BG375").Value = 109
$ws.Range("BG376").Value = 23
$ws.Range("BG377").Value = 57
$ws.Range("BG378").Value = 31
$ws.Range("BG379").Value = 1
$ws.Range("BG380").Value = 7
$ws.Range("BG384").Value = 3
$ws.Range("BG386").Value = 23
$ws.Range("BG388").Value = 26
$ws.Range("BG389").Value = 13
$ws.Range("BG390").Value = 6
$ws.Range("BG391").Value = 21
$ws.Range("BG392").Value = 18
$ws.Range("BG395").Value = 29
$ws.Range("BG398").Value = 8
$ws.Range("BG399").Value = 12
$ws.Range("BG400").Value = 27
$ws.Range("BG401").Value = 30
$ws.Range("BG402").Value = 34

# Column width for column BG (59) - approximate bestFit width
$ws.Columns.Item(59).ColumnWidth = 18.83

# Sheet view: scroll window so row 382 / column AJ area is visible, then
# select the full AQ column (matches the target selection range).
$win = $excel.ActiveWindow
$win.ScrollRow = 382
$win.ScrollColumn = 36
$ws.Range("AQ1:AQ1048576").Select()
